# Refresh the daily crypto snapshot: updated Price (col D) / Volume(1h)
# (col E) figures for most rows, and three rows where the underlying API
# reordered two adjacent coins (VeChain/ImmutableX and
# MXToken/RocketPoolETH), swapping their Coin/Link/Price/Volume cells.
#
# Price-column values that look numeric (e.g. "215.38") are written with
# a leading apostrophe - exactly as typing them into the Excel UI would -
# so they stay literal text instead of being auto-converted to numbers
# (the source keeps every Price cell as text, e.g. "27.577.36" isn't a
# valid number anyway). The Style reset afterwards clears the transient
# "quote prefix" cell style that entry mode implies, since these cells
# carried no explicit style before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.605.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '''1.665.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''215.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''23.61'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.29%  '
$ws.Range('D9').Value = '''0.263'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('D12').Value = '''1.901.52'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('D13').Value = '''1.693.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').Value = '''4.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '''0.557'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').Value = '''66.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '''247.47'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('D18').Value = '''27.608.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').Value = '''0.0₃0731'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.63%  '
$ws.Range('D20').Value = '''7.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.63%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  -3.69%  '
$ws.Range('E23').Value = '  -5.10%  '
$ws.Range('E24').Value = '  -4.93%  '
$ws.Range('D25').Value = '''146.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.01%  '
$ws.Range('D26').Value = '''7.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.85%  '
$ws.Range('D27').Value = '''16.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '''0.112'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.72%  '
$ws.Range('E30').Value = '  +3.73%  '
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('D32').Value = '''3.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').Value = '''1.478.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').Value = '''3.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.27%  '
$ws.Range('E35').Value = '  -5.98%  '
$ws.Range('D36').Value = '''0.937'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''0.573'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.20%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0172'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('D40').Value = '''69.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.77%  '
$ws.Range('E41').Value = '  -5.72%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '''5.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.30%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '''1.809.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.50%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '''2.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('E47').Value = '  -4.13%  '
$ws.Range('D48').Value = '''89.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('E49').Value = '  -2.22%  '
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('D51').Value = '''7.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.75%  '
